# Insert two new data rows right before the current row 1155, shifting the
# existing rows 1155-1250 down to 1157-1252 (dimension grows from
# A1:R1250 to A1:R1252).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1155:1156").Insert()

# New row 1155: Lechuga / Conconina(o) / Segunda
$ws.Range("A1155").Value2 = 10
$ws.Range("B1155").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1155").Value2 = "La Araucanía"
$ws.Range("D1155").Value2 = 44769
$ws.Range("E1155").Value2 = 9
$ws.Range("F1155").Value2 = 100112033
$ws.Range("G1155").Value2 = "Lechuga"
$ws.Range("H1155").Value2 = "Conconina(o)"
$ws.Range("I1155").Value2 = "Segunda"
$ws.Range("J1155").Value2 = 200
$ws.Range("K1155").Value2 = 6000
$ws.Range("L1155").Value2 = 6000
$ws.Range("M1155").Value2 = 6000
$ws.Range("N1155").Value2 = "`$/caja 10 unidades"
$ws.Range("O1155").Value2 = "Provincia del Elquí"
$ws.Range("P1155").Value2 = 600
$ws.Range("Q1155").Value2 = 10
$ws.Range("R1155").Value2 = "Hortaliza"

# New row 1156: Lechuga / Escarola / Primera
$ws.Range("A1156").Value2 = 10
$ws.Range("B1156").Value2 = "Vega Modelo de Temuco"
$ws.Range("C1156").Value2 = "La Araucanía"
$ws.Range("D1156").Value2 = 44769
$ws.Range("E1156").Value2 = 9
$ws.Range("F1156").Value2 = 100112033
$ws.Range("G1156").Value2 = "Lechuga"
$ws.Range("H1156").Value2 = "Escarola"
$ws.Range("I1156").Value2 = "Primera"
$ws.Range("J1156").Value2 = 600
$ws.Range("K1156").Value2 = 13000
$ws.Range("L1156").Value2 = 13000
$ws.Range("M1156").Value2 = 13000
$ws.Range("N1156").Value2 = "`$/caja 15 unidades"
$ws.Range("O1156").Value2 = "Provincia del Elquí"
$ws.Range("P1156").Value2 = 867
$ws.Range("Q1156").Value2 = 15
$ws.Range("R1156").Value2 = "Hortaliza"
